$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "26.963.11"
Set-TextValue "E2" "  +0.32%  "
Set-TextValue "D3" "1.818.45"
Set-TextValue "E3" "  +0.47%  "
Set-TextValue "E4" "  +0.16%  "
Set-TextValue "D5" "309.91"
Set-TextValue "E5" "  +0.06%  "
Set-TextValue "E6" "  +0.11%  "
Set-TextValue "D7" "0.4651"
Set-TextValue "D8" "0.3660"
Set-TextValue "E8" "  -1.12%  "
Set-TextValue "D9" "0.07352"
Set-TextValue "E10" "  -0.50%  "
Set-TextValue "E11" "  -1.05%  "
Set-TextValue "D12" "1.823.80"
Set-TextValue "E12" "  -0.75%  "
Set-TextValue "D13" "5.398"
Set-TextValue "E13" "  +0.80%  "
Set-TextValue "D14" "0.07113"
Set-TextValue "E14" "  +1.04%  "
Set-TextValue "E15" "  -0.09%  "
Set-TextValue "D16" "91.47"
Set-TextValue "E16" "  -0.29%  "
Set-TextValue "E17" "  +0.17%  "
Set-TextValue "D18" "0.000008707"
Set-TextValue "E18" "  +0.19%  "
Set-TextValue "E19" "  +0.11%  "
Set-TextValue "E20" "  -0.62%  "
Set-TextValue "D21" "26.979.25"
Set-TextValue "E21" "  +0.31%  "
Set-TextValue "D22" "5.297"
Set-TextValue "E22" "  -0.41%  "
Set-TextValue "D23" "10.57"
Set-TextValue "E23" "  -0.08%  "
Set-TextValue "D24" "2.046.19"
Set-TextValue "E24" "  -0.24%  "
Set-TextValue "D25" "1.895"
Set-TextValue "D26" "150.55"
Set-TextValue "E26" "  -0.66%  "
Set-TextValue "D27" "18.37"
Set-TextValue "E27" "  -0.09%  "
Set-TextValue "D28" "2.139"
Set-TextValue "E28" "  -0.75%  "
Set-TextValue "D29" "5.248"
Set-TextValue "E29" "  -1.47%  "
Set-TextValue "D30" "116.36"
Set-TextValue "E30" "  +0.45%  "
Set-TextValue "D31" "0.08896"
Set-TextValue "E31" "  -0.17%  "
Set-TextValue "D32" "0.7590"
Set-TextValue "E32" "  +0.38%  "
Set-TextValue "D33" "1.165"
Set-TextValue "E33" "  +0.54%  "
Set-TextValue "D34" "4.505"
Set-TextValue "E34" "  +1.06%  "
Set-TextValue "D35" "2.902"
Set-TextValue "E35" "  -0.68%  "
Set-TextValue "E36" "  +0.14%  "
Set-TextValue "D37" "1.088"
Set-TextValue "E37" "  -1.21%  "
Set-TextValue "E38" "  +0.96%  "
Set-TextValue "D39" "0.01946"
Set-TextValue "E39" "  -1.14%  "
Set-TextValue "D40" "2.970"
Set-TextValue "E40" "  +1.40%  "
Set-TextValue "D41" "7.176"
Set-TextValue "E41" "  -0.11%  "
Set-TextValue "D42" "0.5281"
Set-TextValue "E42" "  -0.84%  "
Set-TextValue "D43" "2.351"
Set-TextValue "E43" "  -3.39%  "
Set-TextValue "E44" "  -0.34%  "
Set-TextValue "D45" "8.445"
Set-TextValue "E45" "  -0.50%  "
Set-TextValue "D46" "0.4868"
Set-TextValue "E46" "  -2.20%  "
Set-TextValue "D47" "10.47"
Set-TextValue "E47" "  +1.66%  "
Set-TextValue "E48" "  +0.13%  "
Set-TextValue "E49" "  -0.15%  "
Set-TextValue "D50" "103.50"
Set-TextValue "E50" "  -0.40%  "
Set-TextValue "D51" "0.06294"
Set-TextValue "E51" "  -0.02%  "
